$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reprioritize the todo list: "Specify download folder" (the feature this
# commit implements download-folder selection for) moves to top priority,
# several previously-unprioritized tasks receive explicit priorities, and
# two new tasks are appended at the bottom of the backlog.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 'Specify download folder'
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 'Make sure one can start more than 2 downloads'
$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 'Bug: cancelled download changes to "starting" when new download is started'
$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 'Show "loading'' while data is loading'
$ws.Range("A6").Value = 60
$ws.Range("B6").Value = 'Delete downloaded file after download canceled'
$ws.Range("A7").Value = 10000
$ws.Range("B7").Value = 'Blob size abbreviation'
$ws.Range("A8").Value = 10000
$ws.Range("B8").Value = 'Sort blob list'
$ws.Range("A9").Value = 10000
$ws.Range("B9").Value = 'Sort account list'
$ws.Range("A10").Value = 10000
$ws.Range("B10").Value = 'Add "Refresh" button for accounts'
$ws.Range("A11").Value = 10000
$ws.Range("B11").Value = 'Add "Refresh" button for blobs'
$ws.Range("A12").Value = 10000
$ws.Range("B12").Value = 'Add ability to cancel downloads'
$ws.Range("A13").Value = 10000
$ws.Range("B13").Value = 'Show blob size and other columns'
$ws.Range("A14").Value = 10000
$ws.Range("B14").Value = 'Show detailed error text when blob download fails'
$ws.Range("A15").Value = 10000
$ws.Range("B15").Value = 'Blob download'
$ws.Range("A16").Value = 10000
$ws.Range("B16").Value = 'Add custom accounts'
$ws.Range("A17").Value = 10000
$ws.Range("B17").Value = 'Sort blobs by different fields'
$ws.Range("B18").Value = 'Automatically update account list after account has been added'
$ws.Range("B19").Value = 'Implement remove account'
$ws.Range("B20").Value = 'Set proxy settings'
$ws.Range("B21").Value = 'Show message when there are no items in container'
$ws.Range("B22").Value = 'Protect files from overwriting when downloaded'
$ws.Range("B23").Value = 'Support virtual folders inside blobs'
$ws.Range("B24").Value = 'Add icons for blobs and accounts'
$ws.Range("B25").Value = 'Implement "test access" for accounts'
$ws.Range("B26").Value = 'Add ability to see Blob URL and other data'
$ws.Range("B27").Value = 'Add application icon'
$ws.Range("B28").Value = 'Implement remove download'
$ws.Range("B29").Value = 'Show error when container list download fails'
$ws.Range("B30").Value = 'Show error when blob list download files'
$ws.Range("B31").Value = 'Older downloads should shift to the right'
$ws.Range("B32").Value = 'Open download folder in explorer'
$ws.Range("B33").Value = 'Change folder list background to azure'
$ws.Range("B34").Value = 'Fast multi-threaded download'

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("B6").Select()
